# Updated cryptos list with new Price (D) and Volume(1h) (E) values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.199.67'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.01%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.080.93'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.70%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.31%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '338.94'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.52%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.004'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.24%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5278'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.63%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4366'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.83%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '54.81'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.06%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.09335'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.05%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.172'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.69%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '24.46'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.68%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '8.474'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.75%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.089.88'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.55%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.852'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.12%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '100.56'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.72%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001159'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.33%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.005'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.22%  '
$ws.Range("E19").Value = '  -2.49%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.06715'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.54%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.312'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.26%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.005'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.15%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '30.212.66'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.98%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.40'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.42%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '21.74'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.58%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.840'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +5.52%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '162.13'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.12%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.486'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.89%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '133.43'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.44%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.127'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.37%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.661'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -7.36%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.1048'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.69%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.241'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.04%  '
$ws.Range("E35").Value = '  -1.56%  '
$ws.Range("E36").Value = '  +0.16%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '9.890'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -8.64%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06714'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.74%  '
$ws.Range("E39").Value = '  -1.03%  '
$ws.Range("E40").Value = '  -1.20%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.338'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.22%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.2201'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.35%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6724'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.20%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.367'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.87%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '14.26'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.33%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.004'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.20%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.298'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.75%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.628'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.30%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.212'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.67%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00000000343'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.42%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.209'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.09%  '
